# feat: add 2022-Q3 data
#
# The workbook has a "总计" (totals) summary sheet followed by one sheet per
# quarter (most-recent quarter first). This change adds a new "2022-Q3"
# quarter sheet (with its fund-holding detail rows) right after "总计", and
# records the corresponding summary row at the top of the "总计" data.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. right before
#    the sheet that is currently "2022-Q2").
# ------------------------------------------------------------------
$insertBefore = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($insertBefore)
$newSheet.Name = "2022-Q3"

# Match the page margins used by the rest of the workbook's data sheets
# (0.75in/0.75in/1in/1in, 0.5in header/footer).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Template sheet to copy cell formatting from (same layout/columns).
$template = $wb.Worksheets.Item("2022-Q2")

# Copy the header-row formatting (B1:H1) and the index-column formatting
# (A2:A4) cell-by-cell so the new sheet matches the look of the other
# quarter sheets.
"B1","C1","D1","E1","F1","G1","H1" | ForEach-Object {
    $template.Range($_).Copy()
    $newSheet.Range($_).PasteSpecial(-4122)
}
"A2","A3","A4" | ForEach-Object {
    $template.Range("A2").Copy()
    $newSheet.Range($_).PasteSpecial(-4122)
}

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows. Columns D/E/F/G hold numeric-looking figures that are
# stored as text in this data set, so force a text format before writing
# them (otherwise they'd be parsed as numbers and lose formatting/leading
# zeros).
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "090019"
$newSheet.Range("C2").Value = "大成景恒混合A"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.13"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.98"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.97"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0223"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "001068"
$newSheet.Range("C3").Value = "华融新锐灵活配置混合"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.21"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "53.47"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "6.19"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0130"
$newSheet.Range("H3").Value = 2

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B4").Value = "006038"
$newSheet.Range("C4").Value = "大成景恒混合C"
$newSheet.Range("D4").NumberFormat = "@"
$newSheet.Range("D4").Value = "0.45"
$newSheet.Range("E4").NumberFormat = "@"
$newSheet.Range("E4").Value = "93.98"
$newSheet.Range("F4").NumberFormat = "@"
$newSheet.Range("F4").Value = "1.97"
$newSheet.Range("G4").NumberFormat = "@"
$newSheet.Range("G4").Value = "0.0089"
$newSheet.Range("H4").Value = 5

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: add the 2022-Q3 summary row at the top of
#    the data and shift the existing quarters down by one row.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 0.04

$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 2
$totals.Range("D3").Value = 0.01

$totals.Range("B4").Value = "2021-Q1"
$totals.Range("C4").Value = 3
$totals.Range("D4").Value = 1.16

# New row 5 ("2020-Q4") — copy the index-column format from row 4 first.
$totals.Range("A4").Copy()
$totals.Range("A5").PasteSpecial(-4122)
$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2020-Q4"
$totals.Range("C5").Value = 1
$totals.Range("D5").Value = 1.42
